# Insert a new weekly record at row 44 for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Acelga".
#
# This pushes the existing data rows 44-129 down to 45-130 (dimension grows
# from A1:R129 to A1:R130) and fills the newly opened row 44 with the new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 44:129 down by one row to make room for the new record.
$ws.Rows("44").Insert()

# Populate the new row 44 with the new weekly observation.
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44544
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = 100112009
$ws.Range("G44").Value = "Acelga"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 50
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = 10000
$ws.Range("N44").Value = "$/docena de atados (12 kilos)"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 833
$ws.Range("Q44").Value = 12
$ws.Range("R44").Value = "Hortaliza"
